$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row labeled "grandes regiões e unidades da federação" (row 6) is a
# header-only row with no data. It is removed entirely, and every row below
# it shifts up by one (the data for "norte", "rondônia", etc. all move up).
$ws.Rows.Item(6).Delete()
